$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to D and E columns being updated so Excel
# does not auto-convert numeric-looking strings (e.g. "0.999") into numbers,
# preserving them as text exactly like the source data.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.918.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.542.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.81"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.54"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.542.01"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.484"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.140"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.22%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.427"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000220"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.139.85"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.543.41"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.728.84"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.40"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.65"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.624"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.12"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000130"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.681.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.55"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.96%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +6.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.78"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.529.68"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.07"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "174.72"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0899"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.888"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.64"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.21"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.68"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.29"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.993"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.42%  "
